$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.086.55"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "3.063.13"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.13"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.92"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.060.30"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +5.84%  "
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.69"
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "3.567.89"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "66.059.71"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "3.061.66"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.18"
$ws.Range("E20").Value = "  +16.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "463.46"
$ws.Range("E21").Value = "  +2.71%  "
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.96"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.08"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.21"
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.86"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.47"
$ws.Range("E38").Value = "  +8.09%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.02"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.70"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0361"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "379.44"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "2.756.33"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.10"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.36"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("E51").Value = "  +3.69%  "
